$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "[PyTorch] numpy로부터 텐서 변환(copying과 sharing의 차이)"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/pytorch-tutorial-01"

$ws.Range("D36").Value = "Score-based Generative Models and Diffusion Models"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/352"

$ws.Range("D51").Value = "[git] 브랜치 생성, 변경, 삭제하기"
$ws.Range("E51").Value = "https://bskyvision.com/1241"

$wb.Save()
